# calculator_smell.xlsx — relabel the "smell" questionnaire rows and switch
# the result cell to show a rounded percentage instead of the raw ratio.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two question labels (A2/A3) and the shared "scale 1-10" helper
# text in C2/C3 to the longer, clearer wording. Order matters: the new
# shared-string table entries must be appended in this sequence (A2, A3,
# then C2/C3) to land at the same slots the real workbook ended up with.
$ws.Range("A2").Value = "Smell before illness"
$ws.Range("A3").Value = "Smell during illness"
$ws.Range("C2").Value = "Scale 1-10 (1- no sense of smell, 10 - excellent sence of smell)"
$ws.Range("C3").Value = "Scale 1-10 (1- no sense of smell, 10 - excellent sence of smell)"

# D6 used to just surface the raw probability; now it renders a rounded
# whole-number percentage via CONCAT/ROUND.
$ws.Range("D6").Formula = '=IF(D5<0.5,"Can not be determined",_xlfn.CONCAT(ROUND(D5*100,0), "%"))'

# Columns A and C widen (bestFit) to accommodate the longer label text.
$ws.Columns.Item(1).ColumnWidth = 16.73
$ws.Columns.Item(3).ColumnWidth = 49.02

# Active cell moved to F4 before the file was last saved.
$ws.Range("F4").Select()
